$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '72.972.51'
$ws.Range('E2').Value = '  +2.92%  '

$ws.Range('D3').Value = '3.990.44'
$ws.Range('E3').Value = '  +1.01%  '

$ws.Range('E4').Value = '  -0.07%  '

Set-TextValue $ws.Range('D5') '592.36'
$ws.Range('E5').Value = '  +9.79%  '

Set-TextValue $ws.Range('D6') '159.36'
$ws.Range('E6').Value = '  +7.17%  '

Set-TextValue $ws.Range('D7') '0.687'
$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('E8').Value = '  -0.05%  '

Set-TextValue $ws.Range('D9') '0.750'
$ws.Range('E9').Value = '  +1.63%  '

$ws.Range('E10').Value = '  +1.90%  '

Set-TextValue $ws.Range('D11') '53.85'
$ws.Range('E11').Value = '  -2.91%  '

$ws.Range('E12').Value = '  +0.80%  '

Set-TextValue $ws.Range('D13') '10.95'
$ws.Range('E13').Value = '  +3.44%  '

$ws.Range('D14').Value = '4.620.08'
$ws.Range('E14').Value = '  +0.97%  '

$ws.Range('D15').Value = '3.986.14'
$ws.Range('E15').Value = '  +0.91%  '

$ws.Range('E16').Value = '  +9.60%  '

Set-TextValue $ws.Range('D17') '14.17'
$ws.Range('E17').Value = '  +3.03%  '

Set-TextValue $ws.Range('D18') '20.41'
$ws.Range('E18').Value = '  +0.12%  '

$ws.Range('E19').Value = '  +0.24%  '

$ws.Range('D20').Value = '72.607.20'
$ws.Range('E20').Value = '  +2.59%  '

Set-TextValue $ws.Range('D21') '435.73'
$ws.Range('E21').Value = '  +2.90%  '

Set-TextValue $ws.Range('D22') '4.80'
$ws.Range('E22').Value = '  +14.23%  '

Set-TextValue $ws.Range('D23') '96.16'
$ws.Range('E23').Value = '  -0.83%  '

$ws.Range('E24').Value = '  -4.48%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D25') '4.49'
$ws.Range('E25').Value = '  +20.79%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D26') '14.29'
$ws.Range('E26').Value = '  +0.37%  '

Set-TextValue $ws.Range('D27') '11.24'
$ws.Range('E27').Value = '  -2.07%  '

Set-TextValue $ws.Range('D28') '10.55'
$ws.Range('E28').Value = '  -1.19%  '

Set-TextValue $ws.Range('D29') '5.96'
$ws.Range('E29').Value = '  +2.19%  '

Set-TextValue $ws.Range('D30') '36.38'
$ws.Range('E30').Value = '  -0.10%  '

Set-TextValue $ws.Range('D31') '7.87'
$ws.Range('E31').Value = '  +1.70%  '

Set-TextValue $ws.Range('D32') '13.67'
$ws.Range('E32').Value = '  +2.25%  '

$ws.Range('E33').Value = '  +1.10%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D34') '48.58'
$ws.Range('E34').Value = '  -5.31%  '

$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D35') '676.92'
$ws.Range('E35').Value = '  -1.81%  '

Set-TextValue $ws.Range('D36') '69.64'
$ws.Range('E36').Value = '  +7.60%  '

$ws.Range('D37').Value = '0.0₃0884'
$ws.Range('E37').Value = '  +7.98%  '

Set-TextValue $ws.Range('D38') '0.436'
$ws.Range('E38').Value = '  +0.08%  '

$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D39') '0.147'
$ws.Range('E39').Value = '  -1.98%  '

$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range('D40') '3.37'
$ws.Range('E40').Value = '  -2.35%  '

$ws.Range('E41').Value = '  +0.17%  '

Set-TextValue $ws.Range('D42') '3.35'
$ws.Range('E42').Value = '  +4.56%  '

$ws.Range('E43').Value = '  +0.05%  '

Set-TextValue $ws.Range('D44') '10.84'
$ws.Range('E44').Value = '  +11.26%  '

Set-TextValue $ws.Range('D45') '0.0488'
$ws.Range('E45').Value = '  +1.55%  '

Set-TextValue $ws.Range('D46') '0.150'
$ws.Range('E46').Value = '  +1.21%  '

$ws.Range('E47').Value = '  -2.87%  '

Set-TextValue $ws.Range('D48') '3.41'
$ws.Range('E48').Value = '  +1.00%  '

$ws.Range('E49').Value = '  +1.55%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.810.97'
$ws.Range('E50').Value = '  +11.72%  '

$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D51') '3.39'
$ws.Range('E51').Value = '  +4.45%  '
